$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6 - this shifts the existing rows 6:11 down to 7:12,
# preserving all of their data/formatting (e.g. the date format on column D).
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly record.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C6").Value = "Los Lagos"
$ws.Range("D6").Value = 44645
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 100112012
$ws.Range("G6").Value = "Espinaca"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 25
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("N6").Value = "$/cuna 10 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1000
$ws.Range("Q6").Value = 10
$ws.Range("R6").Value = "Hortaliza"
